$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 241.72223
$ws.Range("I53").Value = 162.5
$ws.Range("J53").Value = 264.35715
$ws.Range("K53").Value = 162.5
$ws.Range("L53").Value = 264.35715
$ws.Range("M53").Value = 474.5
$ws.Range("N53").Value = -1538.35715

$ws.Range("H76").Value = 6950.2856
$ws.Range("I76").Value = 8058.227
$ws.Range("J76").Value = 5075.3076
$ws.Range("K76").Value = 8058.227
$ws.Range("L76").Value = 5075.3076
$ws.Range("M76").Value = -7743.227
$ws.Range("N76").Value = -5705.3076

$ws.Range("H79").Value = 6950.2856
$ws.Range("I79").Value = 8058.227
$ws.Range("J79").Value = 5075.3076
$ws.Range("K79").Value = 8058.227
$ws.Range("L79").Value = 5075.3076
$ws.Range("M79").Value = -6966.227
$ws.Range("N79").Value = -7259.3076

$ws.Range("H98").Value = 1353.8572
$ws.Range("I98").Value = 1387.25
$ws.Range("J98").Value = 1153.5
$ws.Range("K98").Value = 1387.25
$ws.Range("L98").Value = 1153.5
$ws.Range("M98").Value = 110.75
$ws.Range("N98").Value = -4149.5

$ws.Range("H106").Value = 3237
$ws.Range("I106").Value = 2365.2856
$ws.Range("J106").Value = 4762.5
$ws.Range("K106").Value = 2365.2856
$ws.Range("L106").Value = 4762.5
$ws.Range("M106").Value = -1734.2856
$ws.Range("N106").Value = -6024.5

$ws.Range("H107").Value = 245.71428
$ws.Range("I107").Value = 203.63637
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 203.63637
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1716.36363
$ws.Range("N107").Value = -4240

$ws.Range("H116").Value = 37197.168
$ws.Range("I116").Value = 49595.863
$ws.Range("K116").Value = 49595.863
$ws.Range("M116").Value = -46153.863

$ws.Range("H122").Value = 1353.8572
$ws.Range("I122").Value = 1387.25
$ws.Range("J122").Value = 1153.5
$ws.Range("K122").Value = 4161.75
$ws.Range("L122").Value = 3460.5
$ws.Range("M122").Value = -1711.75
$ws.Range("N122").Value = -8360.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 698.6896400000001
$ws.Range("I2").Value = 663.9231
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 663.9231
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -550.9231
$ws.Range("N2").Value = -1226

$ws.Range("H32").Value = 6996.1445
$ws.Range("I32").Value = 4496.8203
$ws.Range("J32").Value = 45985.6
$ws.Range("K32").Value = 4496.8203
$ws.Range("L32").Value = 45985.6
$ws.Range("M32").Value = -4209.8203
$ws.Range("N32").Value = -46559.6

$ws.Range("H45").Value = 1187.5
$ws.Range("I45").Value = 875
$ws.Range("K45").Value = 875
$ws.Range("M45").Value = -498

$ws.Range("H63").Value = 8360.5
$ws.Range("I63").Value = 9700.625
$ws.Range("K63").Value = 9700.625
$ws.Range("M63").Value = -9014.625

$ws.Range("H66").Value = 8360.5
$ws.Range("I66").Value = 9700.625
$ws.Range("K66").Value = 48503.125
$ws.Range("M66").Value = -45071.125

$ws.Range("H110").Value = 1466.6364
$ws.Range("I110").Value = 1276.1428
$ws.Range("J110").Value = 1800
$ws.Range("K110").Value = 1276.1428
$ws.Range("L110").Value = 1800
$ws.Range("M110").Value = 768.8571999999999
$ws.Range("N110").Value = -5890

$ws.Range("H116").Value = 698.6896400000001
$ws.Range("I116").Value = 663.9231
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 663.9231
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1630.0769
$ws.Range("N116").Value = -5588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 698.6896400000001
$ws.Range("I3").Value = 663.9231
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 663.9231
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -549.9231
$ws.Range("N3").Value = -1228

$ws.Range("H64").Value = 392.875
$ws.Range("I64").Value = 195
$ws.Range("J64").Value = 458.83334
$ws.Range("K64").Value = 195
$ws.Range("L64").Value = 458.83334
$ws.Range("M64").Value = 30
$ws.Range("N64").Value = -908.83334

$ws.Range("H67").Value = 392.875
$ws.Range("I67").Value = 195
$ws.Range("J67").Value = 458.83334
$ws.Range("K67").Value = 195
$ws.Range("L67").Value = 458.83334
$ws.Range("M67").Value = 585
$ws.Range("N67").Value = -2018.83334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0

$ws.Range("H92").Value = 38114.285
$ws.Range("J92").Value = 38114.285
$ws.Range("L92").Value = 38114.285
$ws.Range("N92").Value = -43106.285

$ws.Range("H107").Value = 818.6429000000001
$ws.Range("I107").Value = 280.14285
$ws.Range("J107").Value = 1357.1428
$ws.Range("K107").Value = 280.14285
$ws.Range("L107").Value = 1357.1428
$ws.Range("M107").Value = 1639.85715
$ws.Range("N107").Value = -5197.1428

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 620
$ws.Range("I5").Value = 422.5
$ws.Range("J5").Value = 936
$ws.Range("K5").Value = 1267.5
$ws.Range("L5").Value = 2808
$ws.Range("M5").Value = -1155.5
$ws.Range("N5").Value = -3032

$ws.Range("H123").Value = 4980
$ws.Range("J123").Value = 4980
$ws.Range("L123").Value = 14940
$ws.Range("N123").Value = -19840

$ws.Range("H135").Value = 620
$ws.Range("I135").Value = 422.5
$ws.Range("J135").Value = 936
$ws.Range("K135").Value = 3802.5
$ws.Range("L135").Value = 8424
$ws.Range("M135").Value = -1267.5
$ws.Range("N135").Value = -13494

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 566.1429000000001
$ws.Range("I22").Value = 390.75
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 390.75
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -95.75
$ws.Range("N22").Value = -1390

$ws.Range("H27").Value = 566.1429000000001
$ws.Range("I27").Value = 390.75
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 390.75
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -283.75
$ws.Range("N27").Value = -1014

$ws.Range("H46").Value = 547.4
$ws.Range("I46").Value = 463.125
$ws.Range("J46").Value = 603.5833
$ws.Range("K46").Value = 463.125
$ws.Range("L46").Value = 603.5833
$ws.Range("M46").Value = -275.125
$ws.Range("N46").Value = -979.5833

$ws.Range("H55").Value = 303.23077
$ws.Range("J55").Value = 381.66666
$ws.Range("L55").Value = 381.66666
$ws.Range("M55").Value = -727.66666
